$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# No. Expediente Clinico block (A6 merged A6:F7, G6 merged G6:H7)
$ws.Range("A6").Value = "11111111111  22222222222    "
$ws.Range("G6").Value = "/201761854"

# Fecha de Nacimiento (A9) -- looks like a date, force text entry so it is not
# reinterpreted as a date serial number (matches original plain-text storage)
$ws.Range("A9").Value = "'2017-10-04"

# Edad (E9) -- purely numeric string, force text entry so it stays text
$ws.Range("E9").Value = "'2112"

# Hora de la asistencia medica (D14)
$ws.Range("D14").Value = "Hora: 10:41:45"

# Fecha de la asistencia medica (A15)
$ws.Range("A15").Value = "23/10/2017"
